# Update countries & provincias Spain
# Applies the 14-Abril-2020 16:22 data refresh to the "Pais" sheet:
#  - refreshed case counters for a handful of countries (which also
#    re-shuffles their rank/row position in this totals-sorted table)
#  - bumped the "Datos actualizados..." timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 8 (Alemania): updated counters, same rank ---
$ws.Range("B8").Value = 130434
$ws.Range("C8").Value = 362
$ws.Range("E8").Value = 59014
$ws.Range("G8").Value = 26
$ws.Range("H8").Value = 3220

# --- Row 17 (Brasil): updated counters, same rank ---
$ws.Range("B17").Value = 23955
$ws.Range("C17").Value = 525
$ws.Range("E17").Value = 19615
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = 1361

# --- Rows 49-51: Singapur overtakes Republica Dominicana & Finlandia ---
$ws.Range("A49").Value = "Singapur"
$ws.Range("B49").Value = 3252
$ws.Range("C49").Value = 334
$ws.Range("D49").Value = 611
$ws.Range("E49").Value = 2631
$ws.Range("F49").Value = 29
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 10

$ws.Range("A50").Value = "Republica Dominicana"
$ws.Range("B50").Value = 3167
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 152
$ws.Range("E50").Value = 2838
$ws.Range("F50").Value = 147
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 177

$ws.Range("A51").Value = "Finlandia"
$ws.Range("B51").Value = 3161
$ws.Range("C51").Value = 97
$ws.Range("D51").Value = 300
$ws.Range("E51").Value = 2797
$ws.Range("F51").Value = 74
$ws.Range("G51").Value = 5
$ws.Range("H51").Value = 64

# --- Row 80 (Eslovaquia): updated counters, same rank ---
$ws.Range("D80").Value = 113
$ws.Range("E80").Value = 720

# --- Rows 85-86: Bulgaria overtakes Crucero ---
$ws.Range("A85").Value = "Bulgaria"
$ws.Range("B85").Value = 713
$ws.Range("C85").Value = 28
$ws.Range("D85").Value = 81
$ws.Range("E85").Value = 597
$ws.Range("F85").Value = 36
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 35

$ws.Range("A86").Value = "Crucero"
$ws.Range("B86").Value = 712
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 639
$ws.Range("E86").Value = 61
$ws.Range("F86").Value = 7
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 12

# --- Rows 118-119: Mayotte overtakes Kenia ---
$ws.Range("A118").Value = "Mayotte"
$ws.Range("B118").Value = 217
$ws.Range("C118").Value = 10
$ws.Range("D118").Value = 69
$ws.Range("E118").Value = 145
$ws.Range("F118").Value = 3
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 3

$ws.Range("A119").Value = "Kenia"
$ws.Range("B119").Value = 208
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 40
$ws.Range("E119").Value = 159
$ws.Range("F119").Value = 2
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 9

# --- Timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 16:22"
